$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Purchase 22-23"
$ws2 = $wb.Worksheets.Item(2)   # "Sale 22-23"

# ---------------------------------------------------------------------------
# Sheet "Purchase 22-23": the b22-23MQ428 Print House running total no longer
# totals on row 14; a new row 15 is added for the 07/23-24 invoice instead.
# (Done before the "Sale 22-23" edits below so new shared strings land in
# the same table order as the source edit.)
# ---------------------------------------------------------------------------

$ws1.Range("F14").ClearContents()

$ws1.Range("A14:F14").Copy()
$ws1.Range("A15:F15").PasteSpecial(-4122)   # xlPasteFormats
$ws1.Range("B15").Value = 45020
$ws1.Range("C15").Value = "07/23-24"
$ws1.Range("D15").Value = "Namrata Rubber Product"
$ws1.Range("E15").Value = 2649
$ws1.Range("F15").Formula = "=E13+E14+E15"

# ---------------------------------------------------------------------------
# Sheet "Sale 22-23": remove the Marcfremoit (b22-23MQ209) settled block and
# the Renaldo (b22-23MQ319) settled block, update the latest Putzmeister
# advance rows with the new 2023-24 invoices, and renumber the Sr. No column.
# ---------------------------------------------------------------------------

# Update the two most-recent Putzmeister rows to the new 23-24 invoices.
$ws2.Range("B8").Value = 45020
$ws2.Range("C8").Value = "b23-24MQ101"
$ws2.Range("E8").Value = 107945.2

$ws2.Range("B9").Value = 45020
$ws2.Range("C9").Value = "b23-24MQ102"
$ws2.Range("E9").Value = 290498.3

# Remove the settled Marcfremoit block (rows 11-14, incl. the merged summary
# row) - everything below shifts up by 4 rows.
$ws2.Rows("11:14").Delete()

# Remove the settled Renaldo block (now at rows 17-18 after the shift above).
$ws2.Rows("17:18").Delete()

# Renumber the Sr. No (column A) entries that moved up.
$ws2.Range("A11").Value = 3
$ws2.Range("A13").Value = 4
$ws2.Range("A15").Value = 5
$ws2.Range("A17").Value = 6
$ws2.Range("A19").Value = 7
$ws2.Range("A22").Value = 8
$ws2.Range("A23").ClearContents()

# ---------------------------------------------------------------------------
# Active sheet / selections: "Sale 22-23" becomes the active tab, with A23
# selected there; "Purchase 22-23" keeps D20 selected for when it regains
# focus.
# ---------------------------------------------------------------------------

$ws1.Range("D20").Select()
$ws2.Activate()
$ws2.Range("A23").Select()
